# Atualização arquivos usados em aula.
# Adds a new user row (Daniela Bragança) to the access-list table on
# "Planilha1", mirroring the existing rows (Aluno / Usuário / Senha / IP).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# New row 12: name, username, password (no IP, like several other rows).
$ws.Range("A12").Value = "Daniela Bragança"
$ws.Range("B12").Value = "daniela.braganca"
$ws.Range("C12").Value = "DanBra(äsd86Q!2024>"

# Leave the selection where the author left it when saving.
$ws.Range("C19").Select()
